$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two header labels (target dates changed from 04.06/21-04.06 to 04.04/01-04.04)
$ws.Range("M5").Value = "Целевое значение показов 01-04.04. (15 суток)"
$ws.Range("P5").Value = "Факт показов                                          01-04.04 (15 суток)"

# Clear the stray N15 cell (was holding a leftover value of 1000)
$ws.Range("N15").Clear()

# Update the selected/active cell to match the recorded cursor position
$ws.Range("P13").Select()
